$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell updates that create new shared-string entries ---
# (kept in this specific order so the resulting shared-strings table layout
#  matches the target workbook's append order)
$ws.Range("B13").Value = "Orange Pi Zero3"
$ws.Range("D13").Value = "DDR4 4GB"
$ws.Range("C13").Value = "Allwinner H618, 1.5 GHz"
$ws.Range("C5").Value = "Celeron B800, 1.5 GHz"
$ws.Range("C4").Value = "Intel i5-6200U, 2.8 GHz"
$ws.Range("C2").Value = "Ryzen 5 3600, 3.6GHz"
$ws.Range("B8").Value = "NanoPi NEO3"
$ws.Range("C8").Value = "RockChip RK3328, 1.3 GHz"
$ws.Range("C9").Value = "Snapdragon 660, 2.2 GHz"
$ws.Range("C10").Value = "Exynos 7870, 1.6 GHz"
$ws.Range("C11").Value = "ARM Cortex-A55"
$ws.Range("B12").Value = "TV-Box Vontar"
$ws.Range("C12").Value = "Amlogic S905W2, 1 GHz"
$ws.Range("D7").Value = "DDR4 6GB"
$ws.Range("C7").Value = "Helio G90T, 2 GHz"
$ws.Range("E13").Value = "106ms"
$ws.Range("F13").Value = "990ms"
$ws.Range("G13").Value = "12.3s"
$ws.Range("B14").Value = "Raspberry Pi Zero W"
$ws.Range("C14").Value = "Broadcom BCM2835, 1 GHz"
$ws.Range("D14").Value = "DDR3 512MB"
$ws.Range("B16").Value = "Libre La Frite"
$ws.Range("D16").Value = "DDR4 1GB"
$ws.Range("C16").Value = "Amlogic S805X-AC, 1.4 GHz"
$ws.Range("C6").Value = "Intel E8500 Duo, 3.16 GHz"
$ws.Range("D6").Value = "DDR3 8GB"
$ws.Range("B6").Value = "PC Dell"
$ws.Range("E6").Value = "70ms"
$ws.Range("F6").Value = "700ms"
$ws.Range("G6").Value = "6.5s"
$ws.Range("B15").Value = "Raspberry Pi B"
$ws.Range("C15").Value = "ARM1176JZF-S, 700 МГц"
$ws.Range("D15").Value = "DDR2 256MB"

# --- Cell updates that reuse already-existing shared strings ---
$ws.Range("D3").Value = "DDR4 16GB"
$ws.Range("H6").Value = "-"
$ws.Range("D9").Value = "DDR4 8GB"
$ws.Range("D10").Value = "DDR4 2GB"
$ws.Range("A13").Value = "Linux"
$ws.Range("H13").Value = "-"
$ws.Range("A14").Value = "Linux"
$ws.Range("H14").Value = "-"
$ws.Range("A15").Value = "Linux"
$ws.Range("A16").Value = "Linux"

# --- Adjust column widths (closest achievable given COM width quantization) ---
$ws.Columns.Item(2).ColumnWidth = 22.666666666666668
$ws.Columns.Item(3).ColumnWidth = 25.5
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666

# --- Update selection to match target (D17) ---
$ws.Range("D17").Select()
